$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 4,23
$arr[0,0] = 0.00495867768595041
$arr[0,1] = 0.980165289256198
$arr[0,2] = 0.028099173553719
$arr[0,3] = 0.966942148760331
$arr[0,4] = 0.963636363636364
$arr[0,5] = 0.00495867768595041
$arr[0,6] = 0.00330578512396694
$arr[0,7] = 0.990082644628099
$arr[0,8] = 0.00165289256198347
$arr[0,9] = 0.986776859504132
$arr[0,10] = 0
$arr[0,11] = 0.892561983471074
$arr[0,12] = 0.0479338842975207
$arr[0,13] = 0.056198347107438
$arr[0,14] = 0.998347107438017
$arr[0,15] = 0.00495867768595041
$arr[0,16] = 0
$arr[0,17] = 0.0115702479338843
$arr[0,18] = 0.00661157024793388
$arr[0,19] = 0
$arr[0,20] = 0.975206611570248
$arr[0,21] = 0.00991735537190083
$arr[0,22] = 0.00330578512396694
$arr[1,0] = 0.976859504132231
$arr[1,1] = 0.00495867768595041
$arr[1,2] = 0.00991735537190083
$arr[1,3] = 0
$arr[1,4] = 0.00661157024793388
$arr[1,5] = 0.927272727272727
$arr[1,6] = 0
$arr[1,7] = 0.00495867768595041
$arr[1,8] = 0.985123966942149
$arr[1,9] = 0.00661157024793388
$arr[1,10] = 0.988429752066116
$arr[1,11] = 0.0809917355371901
$arr[1,12] = 0.938842975206612
$arr[1,13] = 0.940495867768595
$arr[1,14] = 0
$arr[1,15] = 0.00330578512396694
$arr[1,16] = 0
$arr[1,17] = 0
$arr[1,18] = 0.988429752066116
$arr[1,19] = 0.991735537190083
$arr[1,20] = 0.00495867768595041
$arr[1,21] = 0
$arr[1,22] = 0.00165289256198347
$arr[2,0] = 0.00991735537190083
$arr[2,1] = 0.0115702479338843
$arr[2,2] = 0.947107438016529
$arr[2,3] = 0.00991735537190083
$arr[2,4] = 0.028099173553719
$arr[2,5] = 0.00165289256198347
$arr[2,6] = 0.00330578512396694
$arr[2,7] = 0.00165289256198347
$arr[2,8] = 0.00495867768595041
$arr[2,9] = 0.00495867768595041
$arr[2,10] = 0.00826446280991736
$arr[2,11] = 0.0231404958677686
$arr[2,12] = 0.00330578512396694
$arr[2,13] = 0.00165289256198347
$arr[2,14] = 0
$arr[2,15] = 0.991735537190083
$arr[2,16] = 1
$arr[2,17] = 0.986776859504132
$arr[2,18] = 0.00495867768595041
$arr[2,19] = 0.00330578512396694
$arr[2,20] = 0.00661157024793388
$arr[2,21] = 0.983471074380165
$arr[2,22] = 0.990082644628099
$arr[3,0] = 0.00826446280991736
$arr[3,1] = 0.00330578512396694
$arr[3,2] = 0.0148760330578512
$arr[3,3] = 0.0214876033057851
$arr[3,4] = 0
$arr[3,5] = 0.0644628099173554
$arr[3,6] = 0.993388429752066
$arr[3,7] = 0.00330578512396694
$arr[3,8] = 0.00826446280991736
$arr[3,9] = 0.00165289256198347
$arr[3,10] = 0.00330578512396694
$arr[3,11] = 0.00330578512396694
$arr[3,12] = 0.00826446280991736
$arr[3,13] = 0.00165289256198347
$arr[3,14] = 0.00165289256198347
$arr[3,15] = 0
$arr[3,16] = 0
$arr[3,17] = 0
$arr[3,18] = 0
$arr[3,19] = 0.00495867768595041
$arr[3,20] = 0.0132231404958678
$arr[3,21] = 0.00661157024793388
$arr[3,22] = 0.00330578512396694
$ws.Range("B2:X5").Value = $arr
Write-Output "done"
